# Insert a new "Charges" parameter row into the Defaults table on Sheet1.
# The table lives in A1:D8 (header in row 1). We need to insert a new row
# right after the "NoiseFilter" row (row 3) and before the "IsotopicPercentage"
# row (old row 4), shifting everything below down by one, then fill the new
# row 4 with the Charges parameter data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 4, shifting rows 4:8 down to 5:9.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the Charges parameter.
$ws.Range("A4").Value = "Charges"
$ws.Range("B4").Value = "1,2"
$ws.Range("C4").Value = "Everytime"
$ws.Range("D4").Value = "The range of charges to test. List charges separated by a comma"

# Update the selection to match the authored state (range A4:D4 selected).
$ws.Range("A4:D4").Select()
